$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A22").Value = 2021
$ws.Range("B22").Value = 64
$ws.Range("C22").Formula = "=53+16+6"
$ws.Range("D22").Value = 14
$ws.Range("E22").Value = 54
$ws.Range("F22").Value = 10
$ws.Range("G22").Value = 0.64
$ws.Range("H22").Value = 0.53
$ws.Range("I22").Value = 0.13
$ws.Range("J22").Formula = "=E22/H22"
$ws.Range("K22").Formula = "=E22*I22"
$ws.Range("L22").Formula = "=F22/H22"
$ws.Range("M22").Formula = "=F22*I22"
$ws.Range("N22").Value = "##not finalized, new data, check with Agnes"

$ws.Range("H23").Select()
